$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 581.92
$ws.Range("I15").Value = 581.92
$ws.Range("K15").Value = 1745.76
$ws.Range("M15").Value = -1576.76
$ws.Range("H40").Value = 1596.6666
$ws.Range("J40").Value = 1775
$ws.Range("L40").Value = 1775
$ws.Range("N40").Value = -2125
$ws.Range("H69").Value = 2200
$ws.Range("I69").Value = 1500
$ws.Range("J69").Value = 2375
$ws.Range("K69").Value = 4500
$ws.Range("L69").Value = 7125
$ws.Range("M69").Value = -3626
$ws.Range("N69").Value = -8873
$ws.Range("H72").Value = 2200
$ws.Range("I72").Value = 1500
$ws.Range("J72").Value = 2375
$ws.Range("K72").Value = 13500
$ws.Range("L72").Value = 21375
$ws.Range("M72").Value = -9132
$ws.Range("N72").Value = -30111
$ws.Range("H98").Value = 6716.25
$ws.Range("I98").Value = 4407.8125
$ws.Range("J98").Value = 15950
$ws.Range("K98").Value = 4407.8125
$ws.Range("L98").Value = 15950
$ws.Range("M98").Value = -2909.8125
$ws.Range("N98").Value = -18946
$ws.Range("H122").Value = 6716.25
$ws.Range("I122").Value = 4407.8125
$ws.Range("J122").Value = 15950
$ws.Range("K122").Value = 13223.4375
$ws.Range("L122").Value = 47850
$ws.Range("M122").Value = -10773.4375
$ws.Range("N122").Value = -52750

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H112").Value = 21995.75
$ws.Range("J112").Value = 21995.75
$ws.Range("L112").Value = 21995.75
$ws.Range("N112").Value = -24949.75
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1302.0555
$ws.Range("I58").Value = 1272.4073
$ws.Range("J58").Value = 1391
$ws.Range("K58").Value = 1272.4073
$ws.Range("L58").Value = 1391
$ws.Range("M58").Value = -1069.4073
$ws.Range("N58").Value = -1797
$ws.Range("H125").Value = 13518.5
$ws.Range("J125").Value = 13518.5
$ws.Range("L125").Value = 13518.5
$ws.Range("N125").Value = -18438.5
$ws.Range("H136").Value = 1302.0555
$ws.Range("I136").Value = 1272.4073
$ws.Range("J136").Value = 1391
$ws.Range("K136").Value = 3817.2219
$ws.Range("L136").Value = 4173
$ws.Range("M136").Value = -1267.2219
$ws.Range("N136").Value = -9273
$ws.Range("H141").Value = 17499.8
$ws.Range("J141").Value = 17499.8
$ws.Range("L141").Value = 17499.8
$ws.Range("N141").Value = -27859.8

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 72.625
$ws.Range("I7").Value = 72.625
$ws.Range("K7").Value = 217.875
$ws.Range("M7").Value = -105.875
$ws.Range("H23").Value = 229.2
$ws.Range("I23").Value = 161.33333
$ws.Range("J23").Value = 274.44446
$ws.Range("K23").Value = 483.99999
$ws.Range("L23").Value = 823.33338
$ws.Range("M23").Value = -248.99999
$ws.Range("N23").Value = -1293.33338
$ws.Range("H33").Value = 1613
$ws.Range("I33").Value = 550
$ws.Range("J33").Value = 2676
$ws.Range("K33").Value = 3300
$ws.Range("L33").Value = 16056
$ws.Range("M33").Value = -3017
$ws.Range("N33").Value = -16622
$ws.Range("H34").Value = 1221.3636
$ws.Range("I34").Value = 168.5
$ws.Range("J34").Value = 1455.3334
$ws.Range("K34").Value = 505.5
$ws.Range("L34").Value = 4366.0002
$ws.Range("M34").Value = -421.5
$ws.Range("N34").Value = -4534.0002
$ws.Range("H39").Value = 4553.222
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4553.222
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 13659.666
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -14247.666
$ws.Range("H55").Value = 7099.25
$ws.Range("J55").Value = 7256.2856
$ws.Range("L55").Value = 21768.8568
$ws.Range("N55").Value = -22122.8568

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2058.889
$ws.Range("I43").Value = 932.8570999999999
$ws.Range("J43").Value = 6000
$ws.Range("K43").Value = 932.8570999999999
$ws.Range("L43").Value = 6000
$ws.Range("M43").Value = -781.8570999999999
$ws.Range("N43").Value = -6302
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H97").Value = 21477.6
$ws.Range("I97").Value = 25324.285
$ws.Range("J97").Value = 1282.5
$ws.Range("K97").Value = 25324.285
$ws.Range("L97").Value = 1282.5
$ws.Range("M97").Value = -24828.285
$ws.Range("N97").Value = -2274.5
$ws.Range("H122").Value = 2594.122
$ws.Range("I122").Value = 2636.2812
$ws.Range("J122").Value = 2444.2222
$ws.Range("K122").Value = 7908.8436
$ws.Range("L122").Value = 7332.6666
$ws.Range("M122").Value = -5458.8436
$ws.Range("N122").Value = -12232.6666
$ws.Range("H132").Value = 2383.258
$ws.Range("I132").Value = 1939.7916
$ws.Range("J132").Value = 3903.7144
$ws.Range("K132").Value = 5819.3748
$ws.Range("L132").Value = 11711.1432
$ws.Range("M132").Value = -3289.3748
$ws.Range("N132").Value = -16771.1432

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2200
$ws.Range("I93").Value = 1500
$ws.Range("J93").Value = 2550
$ws.Range("K93").Value = 1500
$ws.Range("L93").Value = 2550
$ws.Range("M93").Value = -252
$ws.Range("N93").Value = -5046
$ws.Range("H100").Value = 6177.6665
$ws.Range("I100").Value = 8054.125
$ws.Range("K100").Value = 8054.125
$ws.Range("M100").Value = -7513.125
$ws.Range("H110").Value = 23822
$ws.Range("J110").Value = 23822
$ws.Range("L110").Value = 23822
$ws.Range("N110").Value = -32002
$ws.Range("H132").Value = 6328.7827
$ws.Range("I132").Value = 6894.8887
$ws.Range("J132").Value = 4290.8
$ws.Range("K132").Value = 20684.6661
$ws.Range("L132").Value = 12872.4
$ws.Range("M132").Value = -18154.6661
$ws.Range("N132").Value = -17932.4

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 12326.667
$ws.Range("I52").Value = 8490
$ws.Range("K52").Value = 8490
$ws.Range("M52").Value = -8264
$ws.Range("H122").Value = 27782780
$ws.Range("I122").Value = 50002540
$ws.Range("J122").Value = 8077.5
$ws.Range("K122").Value = 150007620
$ws.Range("L122").Value = 24232.5
$ws.Range("M122").Value = -150005170
$ws.Range("N122").Value = -29132.5
